# CHC_2012 worksheet: add season-record columns (Wins / Losses / Ties).
# Every player row on this roster shares the team's 61-101-0 2012 season
# record, so the three new trailing columns (AD:AF) get the same constant
# values for every data row, with a header row matching the existing
# header styling (bold, centered, thin border).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 55

$winsCol   = "AD"
$lossesCol = "AE"
$tiesCol   = "AF"

# Seed the new header cells from the last existing header cell (AC1) so
# they pick up the same style (bold font, thin border, centered
# horizontal/vertical alignment) instead of the workbook's plain default.
$ws.Range("AC1").Copy($ws.Range($winsCol + "1"))
$ws.Range("AC1").Copy($ws.Range($lossesCol + "1"))
$ws.Range("AC1").Copy($ws.Range($tiesCol + "1"))

$ws.Range($winsCol + "1").Value   = "Wins"
$ws.Range($lossesCol + "1").Value = "Losses"
$ws.Range($tiesCol + "1").Value   = "Ties"

# Fill in the season record for every player row.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Range($winsCol + $r).Value   = 61
    $ws.Range($lossesCol + $r).Value = 101
    $ws.Range($tiesCol + $r).Value   = 0
}
